# Adds a new "2021" column (M) of data to the worksheet, mirroring the
# formatting already used in column L, and updates the sheet view so the
# selection sits on the new column and the frozen/scrolled view resets to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column M values, cloning the number/font formatting from column L ---
# Each pair is (target cell, source cell to copy formatting from, new value)
$rows = @(
    @{ Target = "M4";  Source = "L4";  Value = 2021 },
    @{ Target = "M5";  Source = "L5";  Value = 2.6 },
    @{ Target = "M6";  Source = "L6";  Value = 3.4 },
    @{ Target = "M7";  Source = "L7";  Value = 1.4 },
    @{ Target = "M8";  Source = "L8";  Value = 3 },
    @{ Target = "M9";  Source = "L9";  Value = 3.5 },
    @{ Target = "M10"; Source = "L10"; Value = 2.4 },
    @{ Target = "M11"; Source = "L11"; Value = 94.4 },
    @{ Target = "M12"; Source = "L12"; Value = 93.1 },
    @{ Target = "M13"; Source = "L13"; Value = 96.2 }
)

foreach ($row in $rows) {
    $ws.Range($row.Source).Copy()
    $ws.Range($row.Target).PasteSpecial(-4122)
    $ws.Range($row.Target).Value = $row.Value
}

$excel.CutCopyMode = 0

# --- Reset the view: scroll back to column A and move the selection to N1 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N1").Select()
